# Scene.xlsx server-data update ("modified for server's data")
#
# The authored change removes two obsolete scene rows (the old "CloneScene"
# row 2 and the old "RebellerNoob" row 4) and tweaks data on the two rows
# that remain:
#   - the surviving "PioneerNoob / villageScene" row gets a new RelivePos
#   - the surviving "Demo1" row gets its ID changed from 4 to 2
#
# Deleting whole rows (rather than just blanking cells) shifts everything
# up and lets the sheet's used range shrink from A1:K6 to A1:K4, matching
# the target layout, while Excel automatically keeps per-cell styles intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 4 first (RebellerNoob / SelectScene) so row indices
# for the still-to-be-deleted row 2 aren't disturbed.
$ws.Rows("4").Delete()
# Remove the old row 2 (CloneScene / Scene2)
$ws.Rows("2").Delete()

# Former row 3 (PioneerNoob / villageScene) is now row 2: update RelivePos.
$ws.Range("E2").Value = "20,0,60"

# Former row 5 (Demo1) is now row 3: update ID.
$ws.Range("B3").Value = "2"

# Leave the cursor where the author's session ended up.
$ws.Range("F5").Select()
